$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 62
$ws.Range("I2").Value = 236
$ws.Range("J2").Value = 857
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 246
$ws.Range("M2").Value = 11
$ws.Range("N2").Value = 144
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 14
$ws.Range("S2").Value = 96
$ws.Range("T2").Value = 131
$ws.Range("U2").Value = 16
$ws.Range("V2").Value = 1281
$ws.Range("X2").Value = 1352
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 17
$ws.Range("AA2").Value = 6
